{"js": "// Update the date label and every \"a\u00d7b=c\" answer cell in the\n// multiplication table to the new values from the commit.\nconst replacements = [\n  [\"2024-03-03 Sunday\", \"2024-03-04 Monday\"],\n  [\"94\u00d746=4324\", \"40\u00d734=1360\"],\n  [\"78\u00d730=2340\", \"38\u00d749=1862\"],\n  [\"34\u00d741=1394\", \"80\u00d786=6880\"],\n  [\"64\u00d789=5696\", \"52\u00d727=1404\"],\n  [\"12\u00d788=1056\", \"67\u00d775=5025\"],\n  [\"73\u00d755=4015\", \"98\u00d763=6174\"],\n  [\"29\u00d761=1769\", \"56\u00d730=1680\"],\n  [\"41\u00d787=3567\", \"21\u00d790=1890\"],\n  [\"65\u00d760=3900\", \"23\u00d788=2024\"],\n  [\"16\u00d793=1488\", \"58\u00d777=4466\"],\n  [\"76\u00d717=1292\", \"68\u00d782=5576\"],\n  [\"74\u00d792=6808\", \"58\u00d747=2726\"],\n  [\"68\u00d799=6732\", \"24\u00d752=1248\"],\n  [\"89\u00d751=4539\", \"23\u00d755=1265\"],\n  [\"94\u00d735=3290\", \"63\u00d723=1449\"],\n  [\"48\u00d711=528\", \"34\u00d775=2550\"],\n  [\"99\u00d762=6138\", \"89\u00d734=3026\"],\n  [\"18\u00d723=414\", \"67\u00d778=5226\"],\n  [\"80\u00d798=7840\", \"37\u00d792=3404\"],\n  [\"66\u00d748=3168\", \"53\u00d765=3445\"],\n  [\"40\u00d781=3240\", \"62\u00d733=2046\"],\n  [\"30\u00d789=2670\", \"28\u00d745=1260\"],\n  [\"96\u00d755=5280\", \"14\u00d752=728\"],\n  [\"98\u00d721=2058\", \"33\u00d759=1947\"],\n  [\"80\u00d714=1120\", \"69\u00d733=2277\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and every \"a\u00d7b=c\" answer cell in the\n# multiplication table to the new values from the commit.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-03-03 Sunday\", \"2024-03-04 Monday\"),\n  @(\"94\u00d746=4324\", \"40\u00d734=1360\"),\n  @(\"78\u00d730=2340\", \"38\u00d749=1862\"),\n  @(\"34\u00d741=1394\", \"80\u00d786=6880\"),\n  @(\"64\u00d789=5696\", \"52\u00d727=1404\"),\n  @(\"12\u00d788=1056\", \"67\u00d775=5025\"),\n  @(\"73\u00d755=4015\", \"98\u00d763=6174\"),\n  @(\"29\u00d761=1769\", \"56\u00d730=1680\"),\n  @(\"41\u00d787=3567\", \"21\u00d790=1890\"),\n  @(\"65\u00d760=3900\", \"23\u00d788=2024\"),\n  @(\"16\u00d793=1488\", \"58\u00d777=4466\"),\n  @(\"76\u00d717=1292\", \"68\u00d782=5576\"),\n  @(\"74\u00d792=6808\", \"58\u00d747=2726\"),\n  @(\"68\u00d799=6732\", \"24\u00d752=1248\"),\n  @(\"89\u00d751=4539\", \"23\u00d755=1265\"),\n  @(\"94\u00d735=3290\", \"63\u00d723=1449\"),\n  @(\"48\u00d711=528\", \"34\u00d775=2550\"),\n  @(\"99\u00d762=6138\", \"89\u00d734=3026\"),\n  @(\"18\u00d723=414\", \"67\u00d778=5226\"),\n  @(\"80\u00d798=7840\", \"37\u00d792=3404\"),\n  @(\"66\u00d748=3168\", \"53\u00d765=3445\"),\n  @(\"40\u00d781=3240\", \"62\u00d733=2046\"),\n  @(\"30\u00d789=2670\", \"28\u00d745=1260\"),\n  @(\"96\u00d755=5280\", \"14\u00d752=728\"),\n  @(\"98\u00d721=2058\", \"33\u00d759=1947\"),\n  @(\"80\u00d714=1120\", \"69\u00d733=2277\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.MatchWholeWord = $true\n  $find.MatchCase = $true\n  $find.Execute([ref]$null, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$null, [ref]2)\n}\n"}
